$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 3 and 4 (D9 / D10 descriptions)
$ws.Range("B3").Value = "RPWM for actuator 1"
$ws.Range("B4").Value = "LPWM for actuator 1"

# Insert two new rows after row 4 for D6 / D5 actuator 2 entries
$ws.Rows("5:6").Insert()
$ws.Range("A5").Value = "D6"
$ws.Range("B5").Value = "RPWM for actuator 2"
$ws.Range("A6").Value = "D5"
$ws.Range("B6").Value = "LPWM for actuator 2"

$ws.Range("B13").Select()
